# Sync attendance_reports: swap the order of the two "Recorded By" entries
# (column G) for the rows that were re-ordered upstream, e.g.
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"
#   "backup@backdoor.com, System" ->  "System, backup@backdoor.com"
#   "admin@admin.com, dnasr281@gmail.com" -> "dnasr281@gmail.com, admin@admin.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3,4,6,10,11,12,13,14,15,17,18,19,20,21,22,24,26,29,30,32,36,37,38,39,40,41,43,44,45,46,47,48,50,52,55,56,58,62,63,64,65,66,67,69,70,71,72,73,74,76,78,83,84,85,86,87,90,92,93,94,96,99,101,109,110,111,112,113,116,118,119,120,122,125,127,135,136,137,138,139,142,144,145,146,148,151,153)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Text
    $parts = $current.Split(",")
    if ($parts.Length -eq 2) {
        $first = $parts[0].Trim()
        $second = $parts[1].Trim()
        $cell.Value = $second + ", " + $first
    }
}
